# Applies the "updated slides and codes for data type" edits:
#   1) Slide 12 ("Dictionaries"): "key-element pairs" -> "key-value pairs"
#   2) Slide 15 (Hands-On write-up): merge the two runs
#      "Write a program that asks the " + "user…" into a single run
#      "Write a program that asks the user…"

$p = $ppt.ActivePresentation

# --- 1) Slide 12: fix "key-element" -> "key-value" -------------------------
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(2)
$tr12 = $sh12.TextFrame.TextRange

$full12 = $tr12.Text
$oldSentence = "Dictionaries are sets of key-element pairs. "
$newSentence = "Dictionaries are sets of key-value pairs. "

$idx12 = $full12.IndexOf($oldSentence)
if ($idx12 -ge 0) {
    $run12 = $tr12.Characters($idx12 + 1, $oldSentence.Length)
    $run12.Text = $newSentence
}

# --- 2) Slide 15: merge "Write a program that asks the " + "user…" --------
$s15 = $p.Slides.Item(15)
$sh15 = $s15.Shapes.Item(2)
$tr15 = $sh15.TextFrame.TextRange

$full15 = $tr15.Text
$oldPhrase = "Write a program that asks the user"
$newPhrase = "Write a program that asks the user" + [char]8230

$idx15 = $full15.IndexOf($oldPhrase)
if ($idx15 -ge 0) {
    $run15 = $tr15.Characters($idx15 + 1, $oldPhrase.Length + 1)
    $run15.Text = $newPhrase
}
